$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new D,E), shifting old D..K to F..M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats / styles from the old data (now in F:G) into the newly inserted D:E columns
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 38100
$ws.Range("E8").Value = 40600
$ws.Range("H8").Value = 35700
$ws.Range("I8").Value = 39300
$ws.Range("D9").Value = 13100
$ws.Range("E9").Value = 13900
$ws.Range("H9").Value = 11300
$ws.Range("I9").Value = 16400
$ws.Range("D10").Value = 25000
$ws.Range("E10").Value = 26700
$ws.Range("H10").Value = 24400
$ws.Range("I10").Value = 22900
$ws.Range("D12").Value = 9300
$ws.Range("E12").Value = 8500
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 43300
$ws.Range("E17").Value = 43200
$ws.Range("H17").Value = 39600
$ws.Range("I17").Value = 44200
$ws.Range("D18").Value = -5200
$ws.Range("E18").Value = -2600
$ws.Range("H18").Value = -3900
$ws.Range("I18").Value = -4900
$ws.Range("D20").Value = -400
$ws.Range("E20").Value = 400
$ws.Range("D21").Value = -4800
$ws.Range("E21").Value = -1400
$ws.Range("H21").Value = -3000
$ws.Range("I21").Value = -4100
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = 100
$ws.Range("D23").Value = -5800
$ws.Range("E23").Value = -2400
$ws.Range("H23").Value = -3900
$ws.Range("I23").Value = -5000
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -5800
$ws.Range("E26").Value = -2400
$ws.Range("H26").Value = 16900
$ws.Range("I26").Value = -5100
$ws.Range("D27").Value = -5800
$ws.Range("E27").Value = -2400
$ws.Range("H27").Value = 16900
$ws.Range("I27").Value = -5100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 400
$ws.Range("E32").Value = -400
$ws.Range("D33").Value = -5800
$ws.Range("E33").Value = -2400
$ws.Range("H33").Value = -4100
$ws.Range("I33").Value = -5100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -5800
$ws.Range("E35").Value = -2400
$ws.Range("H35").Value = -4100
$ws.Range("I35").Value = -5100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 26000
$ws.Range("E41").Value = 32500
$ws.Range("D42").Value = 66100
$ws.Range("E42").Value = 62300
$ws.Range("D43").Value = 16200
$ws.Range("E43").Value = 14800
$ws.Range("D44").Value = 16100
$ws.Range("E44").Value = 13700
$ws.Range("D45").Value = 6400
$ws.Range("E45").Value = 6400
$ws.Range("D46").Value = 130800
$ws.Range("E46").Value = 129700
$ws.Range("D47").Value = 800
$ws.Range("E47").Value = 1500
$ws.Range("D48").Value = 5900
$ws.Range("E48").Value = 6500
$ws.Range("D49").Value = 500
$ws.Range("E49").Value = 500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 3500
$ws.Range("E52").Value = 3700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 141500
$ws.Range("E54").Value = 142000
$ws.Range("D57").Value = 16100
$ws.Range("E57").Value = 15100
$ws.Range("D58").Value = 20000
$ws.Range("E58").Value = 20000
$ws.Range("D59").Value = 47700
$ws.Range("E59").Value = 46800
$ws.Range("D60").Value = 83900
$ws.Range("E60").Value = 81900
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 40100
$ws.Range("E62").Value = 39300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 123900
$ws.Range("E66").Value = 121100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -265800
$ws.Range("E72").Value = -260000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 17600
$ws.Range("E76").Value = 20800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -5800
$ws.Range("E81").Value = -2400
$ws.Range("H81").Value = -4100
$ws.Range("I81").Value = -5100
$ws.Range("D83").Value = 800
$ws.Range("E83").Value = 800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -1300
$ws.Range("E89").Value = 9100
$ws.Range("D91").Value = -300
$ws.Range("E91").Value = -700
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = -200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -3800
$ws.Range("E94").Value = -10100
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -1300
$ws.Range("E100").Value = -1500
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -6400
$ws.Range("E102").Value = -2500
